# Auto update nse_indices_1 outputs
$wb = $excel.ActiveWorkbook

# --- Index Close sheet: row 7 (date 46003) ---
$ws = $wb.Worksheets.Item("Index Close")
$ws.Range("B7").Value = 26046.94921875
$ws.Range("C7").Value = 68509.703125
$ws.Range("D7").Value = 26565.099609375
$ws.Range("F7").Value = 17276.19921875
$ws.Range("G7").Value = 60283.30078125
$ws.Range("I7").Value = 23726.19921875
$ws.Range("K7").Value = 23726.19921875

# --- MTD % sheet: row 7 (date 46003) ---
$ws = $wb.Worksheets.Item("MTD %")
$ws.Range("B7").Value = 5.77
$ws.Range("C7").Value = 2.76
$ws.Range("D7").Value = 5.29
$ws.Range("F7").Value = 7.75
$ws.Range("G7").Value = 6.08
$ws.Range("I7").Value = 4.46
$ws.Range("K7").Value = 4.46

# --- DoD% sheet: row 6 (date 46006) ---
$ws = $wb.Worksheets.Item("DoD%")
$ws.Range("B6").Value = -0.08
$ws.Range("C6").Value = 0.16
$ws.Range("D6").Value = -0.04
$ws.Range("F6").Value = -0.08
$ws.Range("G6").Value = -0.12
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0

# --- DoD% sheet: row 7 (date 46003) ---
$ws.Range("B7").Value = 0.57
$ws.Range("C7").Value = 0.84
$ws.Range("D7").Value = 0.62
$ws.Range("F7").Value = 1.09
$ws.Range("G7").Value = 1.18
$ws.Range("I7").Value = 0.74
$ws.Range("K7").Value = 0.74

# --- Daily Movers sheet: row 6 (date 46006) ---
$ws = $wb.Worksheets.Item("Daily Movers")
$ws.Range("B6").Value = "Nifty500 Multicap 50:25:25, Nifty Next 50, Nifty500 LargeMidSmall Equal-Cap Weighted"
$ws.Range("C6").Value = "Nifty Midcap 100, Nifty 50, Nifty Midcap 50"

# --- Daily Movers sheet: row 7 (date 46003) ---
$ws.Range("B7").Value = "Nifty Midcap 100, Nifty Midcap 150, Nifty Midcap 50"
$ws.Range("C7").Value = "Nifty 50, Nifty 100, Nifty 200"
